# Auto-generated update of leve-profit calculation columns (H-N) across
# several worksheets, per the scheduled pricing-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 226.13924
$ws.Range("I33").Value = 172.15277
$ws.Range("K33").Value = 172.15277
$ws.Range("M33").Value = 56.84723
$ws.Range("H98").Value = 31088.385
$ws.Range("I98").Value = 1254.1305
$ws.Range("K98").Value = 1254.1305
$ws.Range("M98").Value = 243.8695
$ws.Range("H122").Value = 31088.385
$ws.Range("I122").Value = 1254.1305
$ws.Range("K122").Value = 3762.3915
$ws.Range("M122").Value = -1312.3915
$ws.Range("H131").Value = 3092.25
$ws.Range("I131").Value = 2179.1667
$ws.Range("J131").Value = 3640.1
$ws.Range("K131").Value = 6537.500100000001
$ws.Range("L131").Value = 10920.3
$ws.Range("M131").Value = -1497.500100000001
$ws.Range("N131").Value = -21000.3
$ws.Range("H137").Value = 3902.45
$ws.Range("I137").Value = 1160.75
$ws.Range("J137").Value = 6644.15
$ws.Range("K137").Value = 3482.25
$ws.Range("L137").Value = 19932.45
$ws.Range("M137").Value = -932.25
$ws.Range("N137").Value = -25032.45

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H94").Value = 50000
$ws.Range("J94").Value = 50000
$ws.Range("L94").Value = 50000
$ws.Range("N94").Value = -51802
$ws.Range("H122").Value = 1339.1904
$ws.Range("I122").Value = 1284.9445
$ws.Range("J122").Value = 1664.6666
$ws.Range("K122").Value = 3854.8335
$ws.Range("L122").Value = 4993.9998
$ws.Range("M122").Value = -1404.8335
$ws.Range("N122").Value = -9893.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2585.5715
$ws.Range("I105").Value = 1917.6316
$ws.Range("K105").Value = 1917.6316
$ws.Range("M105").Value = -170.6315999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4210.261
$ws.Range("I31").Value = 1542.9615
$ws.Range("J31").Value = 5328.8066
$ws.Range("K31").Value = 1542.9615
$ws.Range("L31").Value = 5328.8066
$ws.Range("M31").Value = -1247.9615
$ws.Range("N31").Value = -5918.8066
$ws.Range("H34").Value = 4210.261
$ws.Range("I34").Value = 1542.9615
$ws.Range("J34").Value = 5328.8066
$ws.Range("K34").Value = 1542.9615
$ws.Range("L34").Value = 5328.8066
$ws.Range("M34").Value = -1340.9615
$ws.Range("N34").Value = -5732.8066
$ws.Range("H132").Value = 102856.71
$ws.Range("I132").Value = 2114.1428
$ws.Range("J132").Value = 203599.28
$ws.Range("K132").Value = 6342.428400000001
$ws.Range("L132").Value = 610797.84
$ws.Range("M132").Value = -3812.428400000001
$ws.Range("N132").Value = -615857.84
$ws.Range("H134").Value = 468131.53
$ws.Range("I134").Value = 1306.4783
$ws.Range("J134").Value = 2001985.2
$ws.Range("K134").Value = 3919.4349
$ws.Range("L134").Value = 6005955.6
$ws.Range("M134").Value = -1384.4349
$ws.Range("N134").Value = -6011025.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 5587.36
$ws.Range("I5").Value = 7165.933
$ws.Range("J5").Value = 3219.5
$ws.Range("K5").Value = 21497.799
$ws.Range("L5").Value = 9658.5
$ws.Range("M5").Value = -21385.799
$ws.Range("N5").Value = -9882.5
$ws.Range("H49").Value = 1200
$ws.Range("J49").Value = 1200
$ws.Range("L49").Value = 3600
$ws.Range("N49").Value = -3912
$ws.Range("H60").Value = 504.875
$ws.Range("I60").Value = 323.16666
$ws.Range("J60").Value = 1050
$ws.Range("K60").Value = 969.4999799999999
$ws.Range("L60").Value = 3150
$ws.Range("M60").Value = -718.4999799999999
$ws.Range("N60").Value = -3652
$ws.Range("H61").Value = 500
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 500
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 1500
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -1930
$ws.Range("H74").Value = 12184.5
$ws.Range("I74").Value = 2000
$ws.Range("J74").Value = 13639.429
$ws.Range("K74").Value = 6000
$ws.Range("L74").Value = 40918.287
$ws.Range("M74").Value = -4939
$ws.Range("N74").Value = -43040.287
$ws.Range("H75").Value = 3489.9
$ws.Range("I75").Value = 200
$ws.Range("J75").Value = 4312.375
$ws.Range("K75").Value = 600
$ws.Range("L75").Value = 12937.125
$ws.Range("M75").Value = 398
$ws.Range("N75").Value = -14933.125
$ws.Range("H76").Value = 4681
$ws.Range("I76").Value = 2400
$ws.Range("J76").Value = 4719.661
$ws.Range("K76").Value = 7200
$ws.Range("L76").Value = 14158.983
$ws.Range("M76").Value = -6817
$ws.Range("N76").Value = -14924.983
$ws.Range("H77").Value = 12184.5
$ws.Range("I77").Value = 2000
$ws.Range("J77").Value = 13639.429
$ws.Range("K77").Value = 18000
$ws.Range("L77").Value = 122754.861
$ws.Range("M77").Value = -12696
$ws.Range("N77").Value = -133362.861
$ws.Range("H78").Value = 3489.9
$ws.Range("I78").Value = 200
$ws.Range("J78").Value = 4312.375
$ws.Range("K78").Value = 1800
$ws.Range("L78").Value = 38811.375
$ws.Range("M78").Value = 3192
$ws.Range("N78").Value = -48795.375
$ws.Range("H79").Value = 4681
$ws.Range("I79").Value = 2400
$ws.Range("J79").Value = 4719.661
$ws.Range("K79").Value = 7200
$ws.Range("L79").Value = 14158.983
$ws.Range("M79").Value = -5874
$ws.Range("N79").Value = -16810.983
$ws.Range("H113").Value = 5113.391
$ws.Range("J113").Value = 964.8333
$ws.Range("L113").Value = 2894.4999
$ws.Range("N113").Value = -7234.4999
$ws.Range("H122").Value = 9611.916999999999
$ws.Range("J122").Value = 18645
$ws.Range("L122").Value = 167805
$ws.Range("N122").Value = -172705
$ws.Range("H135").Value = 5587.36
$ws.Range("I135").Value = 7165.933
$ws.Range("J135").Value = 3219.5
$ws.Range("K135").Value = 64493.397
$ws.Range("L135").Value = 28975.5
$ws.Range("M135").Value = -61958.397
$ws.Range("N135").Value = -34045.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4608.5713
$ws.Range("I80").Value = 5142.857
$ws.Range("J80").Value = 4341.4287
$ws.Range("K80").Value = 5142.857
$ws.Range("L80").Value = 4341.4287
$ws.Range("M80").Value = -4144.857
$ws.Range("N80").Value = -6337.4287
$ws.Range("H83").Value = 4608.5713
$ws.Range("I83").Value = 5142.857
$ws.Range("J83").Value = 4341.4287
$ws.Range("K83").Value = 25714.285
$ws.Range("L83").Value = 21707.1435
$ws.Range("M83").Value = -20722.285
$ws.Range("N83").Value = -31691.1435
$ws.Range("H102").Value = 2352.7856
$ws.Range("I102").Value = 2311.3635
$ws.Range("J102").Value = 2504.6667
$ws.Range("K102").Value = 2311.3635
$ws.Range("L102").Value = 2504.6667
$ws.Range("M102").Value = -689.3634999999999
$ws.Range("N102").Value = -5748.6667
$ws.Range("H122").Value = 1205.7142
$ws.Range("I122").Value = 1240
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 3720
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -1270
$ws.Range("N122").Value = -7900
$ws.Range("H126").Value = 10427.846
$ws.Range("I126").Value = 15420.25
$ws.Range("J126").Value = 2440
$ws.Range("K126").Value = 46260.75
$ws.Range("L126").Value = 7320
$ws.Range("M126").Value = -43790.75
$ws.Range("N126").Value = -12260

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1156.826
$ws.Range("I61").Value = 1154.8636
$ws.Range("J61").Value = 1200
$ws.Range("K61").Value = 1154.8636
$ws.Range("L61").Value = 1200
$ws.Range("M61").Value = -952.8635999999999
$ws.Range("N61").Value = -1604
$ws.Range("H94").Value = 61989.5
$ws.Range("J94").Value = 61989.5
$ws.Range("L94").Value = 61989.5
$ws.Range("N94").Value = -63341.5
$ws.Range("H113").Value = 1156.826
$ws.Range("I113").Value = 1154.8636
$ws.Range("J113").Value = 1200
$ws.Range("K113").Value = 1154.8636
$ws.Range("L113").Value = 1200
$ws.Range("M113").Value = 1015.1364
$ws.Range("N113").Value = -5540

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1576.9459
$ws.Range("I132").Value = 1091.7
$ws.Range("J132").Value = 3656.5715
$ws.Range("K132").Value = 3275.1
$ws.Range("L132").Value = 10969.7145
$ws.Range("M132").Value = -745.1000000000004
$ws.Range("N132").Value = -16029.7145
